$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.456.78'

$ws.Range("D3").Value = '3.003.09'
$ws.Range("E3").Value = '  -0.63%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '545.58'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.25'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.95%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '2.998.98'
$ws.Range("E8").Value = '  -0.59%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.490'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.32%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.87'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +13.34%  '

$ws.Range("E11").Value = '  -0.08%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.448'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.60%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000221'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.52%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.08'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.74%  '

$ws.Range("D15").Value = '3.479.68'
$ws.Range("E15").Value = '  -0.99%  '

$ws.Range("D16").Value = '62.553.60'
$ws.Range("E16").Value = '  +0.96%  '

$ws.Range("D17").Value = '3.007.82'
$ws.Range("E17").Value = '  -0.70%  '

$ws.Range("E18").Value = '  -1.98%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.58'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '471.79'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.28%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.45'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.43%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.655'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.98%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.20'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.68%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.60'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.74%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.63'
$ws.Range("D25").ClearFormats()

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.73'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.15%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.65'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.93%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.01'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +5.63%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.00%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.52'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.92%  '

$ws.Range("E32").Value = '  -2.10%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.35'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.10%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.60'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.80%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '54.73'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.25%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.85'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.03%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '453.28'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.11%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0812'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.83%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0393'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.38%  '

$ws.Range("D40").Value = '2.963.50'
$ws.Range("E40").Value = '  -7.98%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.114'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.53%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.09'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.88%  '

$ws.Range("E43").Value = '  +5.35%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '27.01'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.93%  '

$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.05%  '

$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.249'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.55%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.01'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.86%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.109'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.87%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '115.42'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.45%  '

$ws.Range("D50").Value = '0.0₃0495'
$ws.Range("E50").Value = '  -0.01%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.02'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.12%  '
